{"js": "// Replace the 7 phishing-message paragraphs with their new content.\n// Paragraph indices (0-based, within context.document.body.paragraphs)\n// map 1:1 before/after the edit; only the message body text changes,\n// the surrounding heading/answer paragraphs are untouched.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst replacements = {\n  7: \"Dear John Land,\\u000b \\u000b We've noticed some unusual activity on your account. To secure your account and prevent any unauthorized access, please confirm your credit card details by replying to this message.\\u000b \\u000b We apologize for the inconvenience and appreciate your prompt attention to this matter.\\u000b \\u000b Sincerely,\\u000b [Company Name] Customer Support\",\n  10: \"Daniel Frances, \\u000b\\u000bWe have a host of musicians and bands playing in Malaysia in the upcoming year. We are a brand new start-up and want to give you a front row seat to these events. Buy membership today and recieve your yearly pass to these events using the link below: \\u000b\\u000bTicketmaster.\",\n  13: \"Subject: Important Account Update Required!\\u000b\\u000bDear Barbara Whitaker,\\u000b\\u000bWe hope this message finds you well. We are writing to inform you that your Amazon account requires an immediate update to ensure uninterrupted service.\\u000b\\u000bDue to recent security enhancements, we need to verify your payment details. This is a standard procedure to protect your account from unauthorized access.\\u000b\\u000bPlease reply to this message with the following information:\\u000b\\u000b1. Full Name:\\u000b2. Credit Card Number:\\u000b3. Expiration Date:\\u000b4. CVV Code:\\u000b\\u000bYour prompt response is appreciated. Failure to update your information may result in temporary suspension of your account.\\u000b\\u000bThank you for your cooperation.\\u000b\\u000bBest regards,\\u000b\\u000bAmazon Customer Support\\u000b\\u000bP.S. Please do not reply if you have already updated your information. This is an automated message.\",\n  17: \"Hello Rick,\\u000b\\u000bWe're having some trouble collecting your latest Netflix payment.\\u000b\\u000bPlease click the following link to process your payment: \",\n  19: \"Subject: Important Update from Megabank\\u000b \\u000b Dear Rick,\\u000b \\u000b I hope this message finds you well. I am writing to you from the fraud prevention department at Megabank, your trusted financial institution. \\u000b \\u000b We have detected some suspicious activity on your account and need to verify your identity and credit card information to ensure the security of your funds. Please reply to this message with the following details:\\u000b \\u000b - Full name\\u000b - Date of birth \\u000b - Credit card number\\u000b - Expiration date\\u000b - CVV code\\u000b \\u000b This is a time-sensitive matter, so please respond as soon as possible to avoid any disruption to your account. Your continued trust and security is our top priority.\\u000b \\u000b Thank you for your prompt attention to this important issue.\\u000b \\u000b Sincerely,\\u000b Megabank Fraud Prevention Team\",\n  24: \"Hi Stevie,\\u000b\\u000bA mutual friend mentioned you were looking for a new gym. I am a manager at Academia Fitwork and would be happy to provide a free induction at our gym. \\u000b\\u000bAlternatively, if you would like to sign up for a reduced fee of \\u00a35, please send this to the below details:\\u000b... \\u000b\\u000bThanks, \\u000bJoe\",\n  26: \"Subject: Important Update from HealthWellness.com\\u000b \\u000b Dear Stevie,\\u000b \\u000b I hope this message finds you well. I'm writing to you from HealthWellness.com, a leading provider of health and wellness products and services. We have an important update regarding your account.\\u000b \\u000b Recently, we've noticed some suspicious activity on your HealthWellness.com account. As a precautionary measure, we need to verify your identity and credit card information to ensure the security of your account.\\u000b \\u000b Please reply to this message with the following details:\\u000b \\u000b - Full Name: Stevie Crawford\\u000b - Date of Birth: 1982-03-02\\u000b - Credit Card Number:\\u000b - Expiration Date:\\u000b - CVV Code:\\u000b \\u000b Once we have verified your information, we can secure your account and prevent any unauthorized access. Your privacy and security are our top priorities.\\u000b \\u000b Thank you for your prompt attention to this matter.\\u000b \\u000b Best regards,\\u000b HealthWellness.com Customer Support\",\n};\n\nfor (const [indexStr, newText] of Object.entries(replacements)) {\n  const index = Number(indexStr);\n  paragraphs.items[index].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the 7 phishing-message paragraphs with their new content.\n# Paragraph indices below are 1-based (Word COM convention) and map\n# 1:1 before/after the edit; only the message body text changes, the\n# surrounding heading/answer paragraphs are left untouched.\n$d = $word.ActiveDocument\n\n$replacements = @{\n    8 = \"Dear John Land,`v `v We've noticed some unusual activity on your account. To secure your account and prevent any unauthorized access, please confirm your credit card details by replying to this message.`v `v We apologize for the inconvenience and appreciate your prompt attention to this matter.`v `v Sincerely,`v [Company Name] Customer Support\"\n    11 = \"Daniel Frances, `v`vWe have a host of musicians and bands playing in Malaysia in the upcoming year. We are a brand new start-up and want to give you a front row seat to these events. Buy membership today and recieve your yearly pass to these events using the link below: `v`vTicketmaster.\"\n    14 = \"Subject: Important Account Update Required!`v`vDear Barbara Whitaker,`v`vWe hope this message finds you well. We are writing to inform you that your Amazon account requires an immediate update to ensure uninterrupted service.`v`vDue to recent security enhancements, we need to verify your payment details. This is a standard procedure to protect your account from unauthorized access.`v`vPlease reply to this message with the following information:`v`v1. Full Name:`v2. Credit Card Number:`v3. Expiration Date:`v4. CVV Code:`v`vYour prompt response is appreciated. Failure to update your information may result in temporary suspension of your account.`v`vThank you for your cooperation.`v`vBest regards,`v`vAmazon Customer Support`v`vP.S. Please do not reply if you have already updated your information. This is an automated message.\"\n    18 = \"Hello Rick,`v`vWe're having some trouble collecting your latest Netflix payment.`v`vPlease click the following link to process your payment: \"\n    20 = \"Subject: Important Update from Megabank`v `v Dear Rick,`v `v I hope this message finds you well. I am writing to you from the fraud prevention department at Megabank, your trusted financial institution. `v `v We have detected some suspicious activity on your account and need to verify your identity and credit card information to ensure the security of your funds. Please reply to this message with the following details:`v `v - Full name`v - Date of birth `v - Credit card number`v - Expiration date`v - CVV code`v `v This is a time-sensitive matter, so please respond as soon as possible to avoid any disruption to your account. Your continued trust and security is our top priority.`v `v Thank you for your prompt attention to this important issue.`v `v Sincerely,`v Megabank Fraud Prevention Team\"\n    25 = \"Hi Stevie,`v`vA mutual friend mentioned you were looking for a new gym. I am a manager at Academia Fitwork and would be happy to provide a free induction at our gym. `v`vAlternatively, if you would like to sign up for a reduced fee of \u00a35, please send this to the below details:`v... `v`vThanks, `vJoe\"\n    27 = \"Subject: Important Update from HealthWellness.com`v `v Dear Stevie,`v `v I hope this message finds you well. I'm writing to you from HealthWellness.com, a leading provider of health and wellness products and services. We have an important update regarding your account.`v `v Recently, we've noticed some suspicious activity on your HealthWellness.com account. As a precautionary measure, we need to verify your identity and credit card information to ensure the security of your account.`v `v Please reply to this message with the following details:`v `v - Full Name: Stevie Crawford`v - Date of Birth: 1982-03-02`v - Credit Card Number:`v - Expiration Date:`v - CVV Code:`v `v Once we have verified your information, we can secure your account and prevent any unauthorized access. Your privacy and security are our top priorities.`v `v Thank you for your prompt attention to this matter.`v `v Best regards,`v HealthWellness.com Customer Support\"\n}\n\nforeach ($index in $replacements.Keys) {\n    $d.Paragraphs($index).Range.Text = $replacements[$index]\n}\n"}
